$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 16855.857
$ws.Range("I51").Value = 18749
$ws.Range("K51").Value = 18749
$ws.Range("M51").Value = -18265
$ws.Range("H80").Value = 42774.418
$ws.Range("I80").Value = 20268.8
$ws.Range("K80").Value = 60806.39999999999
$ws.Range("M80").Value = -59808.39999999999
$ws.Range("H83").Value = 42774.418
$ws.Range("I83").Value = 20268.8
$ws.Range("K83").Value = 182419.2
$ws.Range("M83").Value = -177427.2
$ws.Range("H88").Value = 85922480
$ws.Range("J88").Value = 13976188
$ws.Range("L88").Value = 13976188
$ws.Range("N88").Value = -13977000
$ws.Range("H91").Value = 85922480
$ws.Range("J91").Value = 13976188
$ws.Range("L91").Value = 13976188
$ws.Range("N91").Value = -13978996
$ws.Range("H98").Value = 10149.333
$ws.Range("I98").Value = 10179.4
$ws.Range("K98").Value = 10179.4
$ws.Range("M98").Value = -8681.4
$ws.Range("H100").Value = 1492.8462
$ws.Range("I100").Value = 534
$ws.Range("J100").Value = 2611.5
$ws.Range("K100").Value = 534
$ws.Range("L100").Value = 2611.5
$ws.Range("M100").Value = 7
$ws.Range("N100").Value = -3693.5
$ws.Range("H107").Value = 46876330
$ws.Range("I107").Value = 17858662
$ws.Range("K107").Value = 17858662
$ws.Range("M107").Value = -17856742
$ws.Range("H111").Value = 10421821
$ws.Range("I111").Value = 13894318
$ws.Range("K111").Value = 41682954
$ws.Range("M111").Value = -41679887
$ws.Range("H112").Value = 7575.7407
$ws.Range("I112").Value = 1255.5
$ws.Range("J112").Value = 8674.913
$ws.Range("K112").Value = 3766.5
$ws.Range("L112").Value = 26024.739
$ws.Range("M112").Value = -2658.5
$ws.Range("N112").Value = -28240.739
$ws.Range("H122").Value = 10149.333
$ws.Range("I122").Value = 10179.4
$ws.Range("K122").Value = 30538.2
$ws.Range("M122").Value = -28088.2
$ws.Range("H135").Value = 715078.9
$ws.Range("I135").Value = 833762.9399999999
$ws.Range("K135").Value = 7503866.459999999
$ws.Range("M135").Value = -7501331.459999999
$ws.Range("H137").Value = 2337.3333
$ws.Range("I137").Value = 1722
$ws.Range("J137").Value = 2513.1428
$ws.Range("K137").Value = 5166
$ws.Range("L137").Value = 7539.428400000001
$ws.Range("M137").Value = -2616
$ws.Range("N137").Value = -12639.4284
$ws.Range("H141").Value = 1877.0769
$ws.Range("I141").Value = 1877.0769
$ws.Range("K141").Value = 5631.2307
$ws.Range("M141").Value = -451.2307000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2235740
$ws.Range("I32").Value = 2503549
$ws.Range("K32").Value = 2503549
$ws.Range("M32").Value = -2503262
$ws.Range("H74").Value = 65691.375
$ws.Range("I74").Value = 93143
$ws.Range("K74").Value = 93143
$ws.Range("M74").Value = -92269
$ws.Range("H77").Value = 65691.375
$ws.Range("I77").Value = 93143
$ws.Range("K77").Value = 465715
$ws.Range("M77").Value = -461347
$ws.Range("H92").Value = 53947
$ws.Range("J92").Value = 53947
$ws.Range("L92").Value = 53947
$ws.Range("N92").Value = -58939
$ws.Range("H102").Value = 3898.2144
$ws.Range("I102").Value = 3619.1
$ws.Range("J102").Value = 4596
$ws.Range("K102").Value = 3619.1
$ws.Range("L102").Value = 4596
$ws.Range("M102").Value = -1997.1
$ws.Range("N102").Value = -7840
$ws.Range("H110").Value = 33335632
$ws.Range("I110").Value = 2454.3333
$ws.Range("K110").Value = 2454.3333
$ws.Range("M110").Value = -409.3332999999998
$ws.Range("H122").Value = 4354.36
$ws.Range("I122").Value = 3658.2
$ws.Range("K122").Value = 10974.6
$ws.Range("M122").Value = -8524.599999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 530.35
$ws.Range("I94").Value = 350.7143
$ws.Range("K94").Value = 350.7143
$ws.Range("M94").Value = 100.2857
$ws.Range("H134").Value = 6810.2163
$ws.Range("I134").Value = 3279.9473
$ws.Range("K134").Value = 9839.841899999999
$ws.Range("M134").Value = -7304.841899999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = 0
$ws.Range("H31").Value = 7244.143
$ws.Range("I31").Value = 2847.8635
$ws.Range("K31").Value = 2847.8635
$ws.Range("M31").Value = -2552.8635
$ws.Range("H34").Value = 7244.143
$ws.Range("I34").Value = 2847.8635
$ws.Range("K34").Value = 2847.8635
$ws.Range("M34").Value = -2645.8635
$ws.Range("H62").Value = 15628947
$ws.Range("I62").Value = 25002576
$ws.Range("J62").Value = 6233.3335
$ws.Range("K62").Value = 25002576
$ws.Range("L62").Value = 6233.3335
$ws.Range("M62").Value = -25001952
$ws.Range("N62").Value = -7481.3335
$ws.Range("H65").Value = 15628947
$ws.Range("I65").Value = 25002576
$ws.Range("J65").Value = 6233.3335
$ws.Range("K65").Value = 125012880
$ws.Range("L65").Value = 31166.6675
$ws.Range("M65").Value = -125009760
$ws.Range("N65").Value = -37406.6675
$ws.Range("H69").Value = 30000
$ws.Range("I69").Value = 30000
$ws.Range("K69").Value = 30000
$ws.Range("M69").Value = -29251
$ws.Range("H72").Value = 30000
$ws.Range("I72").Value = 30000
$ws.Range("K72").Value = 90000
$ws.Range("M72").Value = -86256
$ws.Range("H74").Value = 99988.5
$ws.Range("J74").Value = 99988.5
$ws.Range("L74").Value = 99988.5
$ws.Range("N74").Value = -101736.5
$ws.Range("H77").Value = 99988.5
$ws.Range("J77").Value = 99988.5
$ws.Range("L77").Value = 299965.5
$ws.Range("N77").Value = -308701.5
$ws.Range("H99").Value = 5331.077
$ws.Range("I99").Value = 3925
$ws.Range("J99").Value = 6209.875
$ws.Range("K99").Value = 3925
$ws.Range("L99").Value = 6209.875
$ws.Range("M99").Value = -2427
$ws.Range("N99").Value = -9205.875
$ws.Range("H126").Value = 5331.077
$ws.Range("I126").Value = 3925
$ws.Range("J126").Value = 6209.875
$ws.Range("K126").Value = 11775
$ws.Range("L126").Value = 18629.625
$ws.Range("M126").Value = -9305
$ws.Range("N126").Value = -23569.625
$ws.Range("H132").Value = 6384
$ws.Range("I132").Value = 2645
$ws.Range("J132").Value = 8627.4
$ws.Range("K132").Value = 7935
$ws.Range("L132").Value = 25882.2
$ws.Range("M132").Value = -5405
$ws.Range("N132").Value = -30942.2
$ws.Range("H134").Value = 8266.796
$ws.Range("J134").Value = 7712.52
$ws.Range("L134").Value = 23137.56
$ws.Range("N134").Value = -28207.56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 15151699
$ws.Range("I14").Value = 15151699
$ws.Range("K14").Value = 45455097
$ws.Range("M14").Value = -45454924
$ws.Range("H34").Value = 4946.3335
$ws.Range("J34").Value = 6165.3335
$ws.Range("L34").Value = 18496.0005
$ws.Range("N34").Value = -18664.0005
$ws.Range("H107").Value = 842.2727
$ws.Range("J107").Value = 920.75
$ws.Range("L107").Value = 2762.25
$ws.Range("N107").Value = -6602.25
$ws.Range("H113").Value = 2467.75
$ws.Range("J113").Value = 3224.3125
$ws.Range("L113").Value = 9672.9375
$ws.Range("N113").Value = -14012.9375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2209.353
$ws.Range("I102").Value = 2222.375
$ws.Range("K102").Value = 2222.375
$ws.Range("M102").Value = -600.375
$ws.Range("H126").Value = 55558056
$ws.Range("I126").Value = 166669140
$ws.Range("J126").Value = 2516.1667
$ws.Range("K126").Value = 500007420
$ws.Range("L126").Value = 7548.500100000001
$ws.Range("M126").Value = -500004950
$ws.Range("N126").Value = -12488.5001
$ws.Range("H132").Value = 19998.75
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9625.375
$ws.Range("I132").Value = 5818.5454
$ws.Range("J132").Value = 12846.538
$ws.Range("K132").Value = 17455.6362
$ws.Range("L132").Value = 38539.614
$ws.Range("M132").Value = -14925.6362
$ws.Range("N132").Value = -43599.614
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 21916448
$ws.Range("I122").Value = 31502356
$ws.Range("K122").Value = 94507068
$ws.Range("M122").Value = -94504618
$ws.Range("H126").Value = 38466670
$ws.Range("I126").Value = 62504336
$ws.Range("K126").Value = 187513008
$ws.Range("M126").Value = -187510538
$ws.Range("H132").Value = 200000
$ws.Range("I132").Value = 100000
$ws.Range("J132").Value = 250000
$ws.Range("K132").Value = 300000
$ws.Range("L132").Value = 750000
$ws.Range("M132").Value = -297470
$ws.Range("N132").Value = -755060
$ws.Range("H136").Value = 49616
$ws.Range("I136").Value = 2284.5
$ws.Range("K136").Value = 6853.5
$ws.Range("M136").Value = -4303.5
